# Status code and Group targeting fixes
# - Duplicate "All Enemies - PC" into a new "All Enemies - Dead Test" sheet
#   placed right after it, so Group/All attacks can be tested against a
#   target list where an enemy has already been removed ("died").
# - On the new sheet, replace the Eagle entry with a new enemy, "Asigaru".

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("All Enemies - PC")
$source.Copy($null, $source)

$newSheet = $wb.Worksheets.Item($source.Index + 1)
$newSheet.Name = "All Enemies - Dead Test"

# Replace the Eagle row with a new enemy entry, "Asigaru".
$newSheet.Range("B5").Value = "Asigaru"

# Make the new sheet the active tab, matching the editor's last view.
$newSheet.Activate()
$newSheet.Range("B6").Select() | Out-Null
